$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("L4").Value = 24
$ws.Range("M4").Value = 24
$ws.Range("N4").Value = 12
$ws.Range("O4").Value = 4

# Row 6
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2
$ws.Range("N6").Value = 6
$ws.Range("O6").Value = 24
$ws.Range("R6").Value = 24
$ws.Range("S6").Value = 12
$ws.Range("T6").Value = 8
$ws.Range("U6").Value = 6

# Row 8
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2
$ws.Range("N8").Value = 3
$ws.Range("O8").Value = 4

# Row 15
$ws.Range("R15").Value = -1
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 3
$ws.Range("U15").Value = 4
$ws.Range("V15").Value = 5
$ws.Range("W15").Value = 1
$ws.Range("X15").Value = 2

# Row 16
$ws.Range("L16").Value = -1
$ws.Range("M16").Value = 1
$ws.Range("N16").Value = 3
$ws.Range("O16").Value = 4

# Row 17
$ws.Range("R17").Formula = "=R15"
$ws.Range("S17").Formula = "=S15*R17"
$ws.Range("T17:X17").Formula = "=T15*S17"

# Row 18
$ws.Range("L18:M18").Formula = "=L16*M18"
$ws.Range("N18").Formula = "=N16*O18"
$ws.Range("O18").Formula = "=O16"
$ws.Range("R18").Formula = "=R15*S18"
$ws.Range("S18:V18").Formula = "=S15*T18"
$ws.Range("W18").Formula = "=W15*X18"
$ws.Range("X18").Formula = "=X15"

# Row 20
$ws.Range("L20").Formula = "=L16"
$ws.Range("M20").Formula = "=M16*L20"
$ws.Range("N20:O20").Formula = "=N16*M20"
$ws.Range("R20").Formula = '=$R$18/R15'
$ws.Range("S20:X20").Formula = '=$R$18/S15'

# Row heights matching the rest of the sheet (18pt, custom height)
$ws.Rows(4).RowHeight = 18
$ws.Rows(6).RowHeight = 18
$ws.Rows(8).RowHeight = 18
$ws.Rows(15).RowHeight = 18
$ws.Rows(16).RowHeight = 18
$ws.Rows(17).RowHeight = 18
$ws.Rows(18).RowHeight = 18
$ws.Rows(20).RowHeight = 18

$excel.Calculate()

$ws.Range("M25").Select() | Out-Null
